$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "BAV21A0F1C"
$ws.Range("C2").Value = "BALL VALVE FL, FB, LP, FLOATING BALL, API 608, A216 GR.WCB, CL 150, RF, B16.5, A105 + ENP BALL, SS410 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, SPW SS304/GRAPH, API 607, LO"
$ws.Range("D2").Value = "4,00"
$ws.Range("E2").Value = "5,00"

# Row 3
$ws.Range("B3").Value = "BAV21A0F1C"
$ws.Range("C3").Value = "BALL VALVE FL, FB, LP, FLOATING BALL, API 608, A216 GR.WCB, CL 150, RF, B16.5, A105 + ENP BALL, SS410 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, SPW SS304/GRAPH, API 607, LO"
$ws.Range("D3").Value = "6,00"
$ws.Range("E3").Value = "7,00"

# Row 4
$ws.Range("B4").Value = "BAV21A0I1C"
$ws.Range("C4").Value = "BALL VALVE FL, FB, LP, FLOATING BALL, API 608, A216 GR.WCB, CL 150, RF, B16.5, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, SPW SS304/GRAPH, API 607, LO"
$ws.Range("E4").Value = "3,00"

# Row 5
$ws.Range("B5").Value = "BAV21A0I1C"
$ws.Range("C5").Value = "BALL VALVE FL, FB, LP, FLOATING BALL, API 608, A216 GR.WCB, CL 150, RF, B16.5, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, SPW SS304/GRAPH, API 607, LO"
$ws.Range("E5").Value = "19,00"

# Row 6
$ws.Range("B6").Value = "BAV24G0I1C"
$ws.Range("C6").Value = "BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO"
$ws.Range("D6").Value = "0,5"
$ws.Range("E6").Value = "1,00"

# Row 7
$ws.Range("B7").Value = "BAV24G0I1C"
$ws.Range("C7").Value = "BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO"
$ws.Range("D7").Value = "0,75"
$ws.Range("E7").Value = "2,00"

# Row 8
$ws.Range("B8").Value = "BAV24G0I1C"
$ws.Range("C8").Value = "BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO"
$ws.Range("D8").Value = "1,00"
$ws.Range("E8").Value = "3,00"

# Row 9
$ws.Range("B9").Value = "CKV21A0B2B"
$ws.Range("C9").Value = "SWING CHECK VALVE FL, API 594, API 598, A216 GR.WCB, CL 150, INST HORIZ/VERT, RF, B16.5, BOLTED COVER, SPW SS304/GRAPH, RENEWABLE SEATS, TRIM #8"
$ws.Range("D9").Value = "2,00"
$ws.Range("E9").Value = "1,00"

# Row 10
$ws.Range("B10").Value = "CKV21A0B2B"
$ws.Range("C10").Value = "SWING CHECK VALVE FL, API 594, API 598, A216 GR.WCB, CL 150, INST HORIZ/VERT, RF, B16.5, BOLTED COVER, SPW SS304/GRAPH, RENEWABLE SEATS, TRIM #8"
$ws.Range("D10").Value = "10,00"
$ws.Range("E10").Value = "3,00"

# Row 11
$ws.Range("B11").Value = "GAV21A0B2B"
$ws.Range("C11").Value = "GATE VALVE FL, API 600, API 598, A216 GR.WCB, CL 150, RF, B16.5, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, FLEXIBLE WEDGE, STEM OS&Y/RSNRO, HO"
$ws.Range("D11").Value = "2,00"
$ws.Range("E11").Value = "1,00"

# Row 12
$ws.Range("B12").Value = "GAV21A0B2B"
$ws.Range("C12").Value = "GATE VALVE FL, API 600, API 598, A216 GR.WCB, CL 150, RF, B16.5, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, FLEXIBLE WEDGE, STEM OS&Y/RSNRO, HO"
$ws.Range("D12").Value = "8,00"
$ws.Range("E12").Value = "8,00"

# Row 13
$ws.Range("B13").Value = "GAV21A0B2B"
$ws.Range("C13").Value = "GATE VALVE FL, API 600, API 598, A216 GR.WCB, CL 150, RF, B16.5, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, FLEXIBLE WEDGE, STEM OS&Y/RSNRO, HO"
$ws.Range("D13").Value = "10,00"
$ws.Range("E13").Value = "11,00"

# Row 14
$ws.Range("B14").Value = "GAV21A0B2B"
$ws.Range("C14").Value = "GATE VALVE FL, API 600, API 598, A216 GR.WCB, CL 150, RF, B16.5, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, FLEXIBLE WEDGE, STEM OS&Y/RSNRO, GO"
$ws.Range("D14").Value = "14,00"
$ws.Range("E14").Value = "1,00"

# Row 15
$ws.Range("B15").Value = "GLV24F0B2B"
$ws.Range("C15").Value = "GLOBE VALVE SW, API 602, API 598, A105, CL 800, SW, B16.11, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, SWIVEL PLUG DISC, STEM OS&Y/RSRO, HO"
$ws.Range("D15").Value = "0,5"
$ws.Range("E15").Value = "1,00"

# Row 16
$ws.Range("B16").Value = "GLV24F0B2B"
$ws.Range("C16").Value = "GLOBE VALVE SW, API 602, API 598, A105, CL 800, SW, B16.11, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, SWIVEL PLUG DISC, STEM OS&Y/RSRO, HO"
$ws.Range("D16").Value = "0,75"
$ws.Range("E16").Value = "9,00"
